# Auto-génération des classes et des specs
# Adds enum hints (on a manual line break) to the "Format" column of
# several rows, normalizes the "mobility" example to upper case, and
# fills in the previously empty "capacity" example.
#
# Note: Range.Find.Execute() on a sub-range (e.g. a table cell's Range)
# is not scoped to that range in this host -- it searches/replaces
# across the whole story. So cell contents are rewritten directly via
# Range.Text assignment instead; [char]11 (vertical tab) is Word's
# manual-line-break character and serializes to <w:br/>.

$d = $word.ActiveDocument

$lineBreak = [char]11

function Set-FormatEnum {
    param($Cell, $EnumText)
    $Cell.Range.Text = "string" + $lineBreak + $EnumText
}

# --- "Type resource" table (Tables(2)) ---
$tResource = $d.Tables(2)

# row 5 -> "type"     : Format string -> string + (Enum : SMUR, SDIS, TSU, SNP, MSPE, SHIP)
Set-FormatEnum $tResource.Rows(5).Cells(3) "(Enum : SMUR, SDIS, TSU, SNP, MSPE, SHIP)"

# row 6 -> "nature"   : Format string -> string + (Enum : EFFECTOR, BASE)
Set-FormatEnum $tResource.Rows(6).Cells(3) "(Enum : EFFECTOR, BASE)"

# row 7 -> "mobility" : Format string -> string + (Enum : FIX, VEHICLE, HELICOPTER, SHIP )
Set-FormatEnum $tResource.Rows(7).Cells(3) "(Enum : FIX, VEHICLE, HELICOPTER, SHIP )"
# ... and its example "Vehicle" -> "VEHICLE"
$tResource.Rows(7).Cells(6).Range.Text = "VEHICLE"

# row 8 -> "capacity" : Format string -> string + (Enum : EMERGENCY, MEDICAL, PARAMEDICAL, UNKNOWN )
Set-FormatEnum $tResource.Rows(8).Cells(3) "(Enum : EMERGENCY, MEDICAL, PARAMEDICAL, UNKNOWN )"
# ... and its example, which was empty, becomes "EMERGENCY"
$tResource.Rows(8).Cells(6).Range.Text = "EMERGENCY"

# --- "Type contact" table (Tables(3)) ---
$tContact = $d.Tables(3)

# row 2 -> "type"     : Format string -> string + (Enum : PMRADD, PHNADD)
Set-FormatEnum $tContact.Rows(2).Cells(3) "(Enum : PMRADD, PHNADD)"
